$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C2 (empty inline string cell)
$ws.Range("C2").Value = ""

# Update C3:C14 with new transcript text
$ws.Range("C3").Value = 'じゃあ えっと まずは 伊藤 電子書籍買えと普通に実際 手本 あの紙で買うほうかどっちらかの意見を聞いていきたいと思いますま 自分は電子 書式の方を選びますかな髪の方で神栖で屋守です理由はまず実際に 紙の本で買った場合ってお金が勝手なまとめ買いとか 漫画とか色々買ってしまうと何かにお金がかかってしまうので電子書籍の場合は会員登録したら 本は読める'
$ws.Range("C4").Value = '本は読めるか辞めるのでも別に読めたら まこりんで鍵はどっちでもいいかどっち行きたいなと思って電話しました 読みたい本 とかなかったらもうそれはもう2回日本でかしかないけど 基本的に連射式の方である場合は 他の紙となくてできるよ'
$ws.Range("C5").Value = 'ここからしゃべっちゃった時なあまあ 僕も 昔 山上 腹立つけど一人暮らし 初めて かなとなとを置く場所なくて泣く泣くまで行った時何だったんですけど 良い点としてはスマホでも読めるよ 待ってまーす お嫁の利点として記憶に残りやすいそういった工夫ができる仲間 神 またただいま マークが引けるとかですけど遊べて工夫が出来たりするのも 強みでも集中できる'
$ws.Range("C6").Value = '集中できる 読書にそういう端末でやったらそういう SNS とか 通知が来たら変な子を集中が割れたりとかパソコンに向けたりできる子っていう できる しそのみたいとこ 千とかで開いたり そういうのは結構あの電子書籍と本を比較した時に残らないから まあ かさばらないし かさばらないし 邪魔にならないけど外出先でも言っても行ける利便性があると思うけど 実物の神様が怒って 実物であったとしても'
$ws.Range("C7").Value = '間違ったとしても インテリアにも代わりにもできる女とかに入れとくだけでも間に5以上 開いてもメールしても使えるし後は電子書籍で買って4の紙でカウントしたら紙で買ったらふる本屋とかあとはのメルカリとかで見れないって処分する時もお金は打ったら還元で帰ってくるってこと考えると 電子書籍よりかはもうお金は返ってくるから得ること考えるんやったらそっちの方が安く済むこともあるやろう しまもも として残るから貸し借り とかもできる いう人間に この漫画 面白かったからどう見たいそうそうそう'
$ws.Range("C8").Value = 'そそそ スマホとかそういうアプリとかやったら スマホと私って読めるかもしれんけど そのすぐ4台とかそういうこと考えたら 真帆の方がいいかな 見て最悪 いらなくなったら売ればお金返ってくるからでも読み終えたときのその中身を見る のは 多分ですけど その調子 の遊びの生死 紙とか漫画とかに'
$ws.Range("C9").Value = '紙とかまあ 漫画とかなんですけどそのあるあるものは 表紙 めくっ町 しなのやつったらその入るでしょう 滋賀工場たワンピースとかあると思うけど あそこになんか遊びとかがあったりしてるんであれば ワンピース パンダマンでもらってるようだったら書いたりとかするんでそういう そういうの意外な遊びとしても楽しめるってのは神'
$ws.Range("C10").Value = 'ハイスクール4で逆に逆に向けたをするためになるのは あと何分ですか 何か他に意見がある人いますか やっぱり その感じがするなっていう風に'
$ws.Range("C11").Value = '押して言ってたんですけどやっぱカニの方が やっぱり他に読み方がちょっと面白み とかがあってその本読んでる感じ もあるかそっちの方をしたいと思いますし なんで それはまた'
$ws.Range("C12").Value = 'それはまだ近い方がいいかもしれないですね ただ 本当におきました 9で紙の方がいいかもしれんけど間違えて'
$ws.Range("C13").Value = '間違えて切っちゃって 8のページ 見ちゃったりましたな'
$ws.Range("C14").Value = 'ドラゴンボールz 今んとこ版24で電子が1で漢字でなってますけど間違いはないですか 何か他に意見 ある人がいれば'

# Clear C15 (empty inline string cell)
$ws.Range("C15").Value = ""

# Update C16 and C17 with new transcript text
$ws.Range("C16").Value = 'あと1分です 意見がなければ そろそろ まとめ の方に入りたいと思うんですけどどうでしょうか まとめると咳が止まらずに手軽に読めるっていう点で 神だと実際に本を読んで本読んで'
$ws.Range("C17").Value = '楽しめるかあります'

# Remove rows 18-20 entirely (they no longer exist in the target sheet)
$ws.Range("A18:C20").EntireRow.Delete()

